$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("AF2").Value = 9.6
$ws.Range("AK2").Value = 1000
$ws.Range("F2").Value = 1.48
$ws.Range("G2").Value = 1.6
$ws.Range("H2").Value = 6.6
$ws.Range("I2").Value = 11
$ws.Range("J2").Value = 3.75
$ws.Range("K2").Value = 4.8
$ws.Range("L2").Value = 1.4
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 3.4
$ws.Range("O2").Value = 1.32
$ws.Range("P2").Value = 1.83
$ws.Range("Q2").Value = 1.96
$ws.Range("R2").Value = 1.32
$ws.Range("S2").Value = 3.4
$ws.Range("T2").Value = 2.06
$ws.Range("U2").Value = 1.73
$ws.Range("V2").Value = 1.12
$ws.Range("W2").Value = 2.66
$ws.Range("X2").Value = 90

# Row 3
$ws.Range("AB3").Value = 7.6
$ws.Range("H3").Value = 4.7
$ws.Range("I3").Value = 5.9
$ws.Range("U3").Value = 1.83
$ws.Range("V3").Value = 1.21

# Row 4
$ws.Range("G4").Value = 1.84
$ws.Range("J4").Value = 3.4
$ws.Range("K4").Value = 4.4
$ws.Range("N4").Value = 2.26
$ws.Range("P4").Value = 1.39
$ws.Range("Q4").Value = 2.34
$ws.Range("T4").Value = 2.08
$ws.Range("U4").Value = 1.59

# Row 5
$ws.Range("AJ5").Value = 220
$ws.Range("AM5").Value = 120
$ws.Range("AN5").Value = 130
$ws.Range("P5").Value = 2.22
$ws.Range("R5").Value = 1.46
$ws.Range("T5").Value = 1.93

# Row 6
$ws.Range("AA6").Value = 1000
$ws.Range("AB6").Value = 15.5
$ws.Range("AC6").Value = 19
$ws.Range("AG6").Value = 12.5
$ws.Range("AH6").Value = 34
$ws.Range("AN6").Value = 3.05
$ws.Range("H6").Value = 15.5
$ws.Range("I6").Value = 16
$ws.Range("J6").Value = 7.8
$ws.Range("L6").Value = 1.21
$ws.Range("N6").Value = 8.800000000000001
$ws.Range("O6").Value = 1.11
$ws.Range("P6").Value = 3.55
$ws.Range("Q6").Value = 1.37
$ws.Range("S6").Value = 1.94
$ws.Range("U6").Value = 2.04
$ws.Range("V6").Value = 1.06

# Row 7
$ws.Range("AI7").Value = 70
$ws.Range("L7").Value = 1.17
$ws.Range("Q7").Value = 1.46
$ws.Range("R7").Value = 1.7
$ws.Range("S7").Value = 2.12
$ws.Range("U7").Value = 1.56

# Row 8
$ws.Range("R8").Value = 1.2

# Row 9
$ws.Range("AB9").Value = 16
$ws.Range("AE9").Value = 20
$ws.Range("AG9").Value = 14
$ws.Range("AN9").Value = 28
$ws.Range("H9").Value = 2.3
$ws.Range("X9").Value = 18

# Row 10
$ws.Range("N10").Value = 4.4
$ws.Range("Q10").Value = 1.81
$ws.Range("U10").Value = 1.89

# Row 11
$ws.Range("AA11").Value = 32
$ws.Range("AI11").Value = 26
$ws.Range("AJ11").Value = 46
$ws.Range("AK11").Value = 27
$ws.Range("AL11").Value = 30
$ws.Range("AM11").Value = 48
$ws.Range("AN11").Value = 17.5
$ws.Range("AO11").Value = 11.5
$ws.Range("H11").Value = 2.36
$ws.Range("I11").Value = 2.38
$ws.Range("P11").Value = 2.62
$ws.Range("R11").Value = 1.67
$ws.Range("V11").Value = 1.72

# Row 12
$ws.Range("AK12").Value = 19.5
$ws.Range("AM12").Value = 48
$ws.Range("F12").Value = 2.3
$ws.Range("G12").Value = 2.32
$ws.Range("K12").Value = 3.95
$ws.Range("P12").Value = 2.7
$ws.Range("R12").Value = 1.71
$ws.Range("S12").Value = 2.32
$ws.Range("W12").Value = 1.76
$ws.Range("X12").Value = 25
$ws.Range("Y12").Value = 21
$ws.Range("Z12").Value = 27

# Row 13
$ws.Range("AJ13").Value = 27
$ws.Range("AN13").Value = 11.5
$ws.Range("F13").Value = 2.18
$ws.Range("H13").Value = 3.65
$ws.Range("I13").Value = 3.7
$ws.Range("J13").Value = 3.7
$ws.Range("K13").Value = 3.75
$ws.Range("O13").Value = 1.24
$ws.Range("Q13").Value = 1.71
$ws.Range("U13").Value = 2.5
$ws.Range("Y13").Value = 16.5
$ws.Range("Z13").Value = 27

# Row 14
$ws.Range("AJ14").Value = 11.5
$ws.Range("AN14").Value = 2.48
$ws.Range("F14").Value = 1.19
$ws.Range("I14").Value = 18.5
$ws.Range("J14").Value = 9.4
$ws.Range("R14").Value = 2.44

# Row 15
$ws.Range("P15").Value = 2.28
$ws.Range("Q15").Value = 1.73
$ws.Range("S15").Value = 2.82
$ws.Range("T15").Value = 1.63
$ws.Range("Y15").Value = 13

# Row 16
$ws.Range("AB16").Value = 16
$ws.Range("AC16").Value = 9.6
$ws.Range("AF16").Value = 21
$ws.Range("AH16").Value = 14.5
$ws.Range("AI16").Value = 36
$ws.Range("AN16").Value = 15
$ws.Range("F16").Value = 2.38
$ws.Range("G16").Value = 2.62
$ws.Range("H16").Value = 2.88
$ws.Range("I16").Value = 3.15
$ws.Range("J16").Value = 3.5
$ws.Range("L16").Value = 1.31
$ws.Range("M16").Value = 1.04
$ws.Range("N16").Value = 4.9
$ws.Range("P16").Value = 2.32
$ws.Range("Q16").Value = 1.61
$ws.Range("R16").Value = 1.53
$ws.Range("S16").Value = 2.5
$ws.Range("T16").Value = 1.53
$ws.Range("U16").Value = 2.48
$ws.Range("V16").Value = 1.46
$ws.Range("W16").Value = 1.61
$ws.Range("X16").Value = 23
$ws.Range("Y16").Value = 18
